$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3) Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4) Forecast data rows: ds, PO_Forecast, yhat_lower, yhat_upper ---
$data = @(
    @(44934.99999999999,81,-54.3025753139986,221.350523849097),
    @(44955.99999999999,80,-62.46981318566891,213.7072678520938),
    @(44962.99999999999,80,-65.71347691380716,214.3159325205371),
    @(44976.99999999999,80,-53.40729397585196,219.4173300851554),
    @(44990.99999999999,79,-49.80708859136582,219.532635254218),
    @(44997.99999999999,79,-70.52197524772404,212.8226613008328),
    @(45004.99999999999,79,-53.10713598152353,219.4790052940547),
    @(45011.99999999999,78,-71.30750542492733,218.7840149041539),
    @(45018.99999999999,78,-70.465691264353,218.0861947886308),
    @(45025.99999999999,78,-55.6542192531793,217.5412095278384),
    @(45032.99999999999,78,-57.67405359736986,217.4181307373896),
    @(45039.99999999999,78,-70.9790842981938,228.1511328488441),
    @(45046.99999999999,77,-71.16056488900254,215.4429453292274),
    @(45053.99999999999,77,-67.7152887766419,223.1827637176355),
    @(45060.99999999999,77,-70.83041802300018,208.1762953952143),
    @(45088.99999999999,76,-58.81840744923019,214.9078535796492),
    @(45095.99999999999,76,-73.17037691955019,216.9168101939839),
    @(45102.99999999999,76,-70.19977603106814,218.1593329561975),
    @(45130.99999999999,75,-66.23641370593647,219.2149337221066),
    @(45144.99999999999,74,-70.42209246650972,216.3104640563868),
    @(45158.99999999999,74,-65.27277832933336,214.3875702189759),
    @(45263.99999999999,71,-73.68159334975979,208.5727699867827),
    @(45354.99999999999,68,-60.09933737996231,208.6558051188399),
    @(45361.99999999999,68,-68.0052253494465,203.7609075901971),
    @(45368.99999999999,68,-60.12858001290352,210.9350216218293),
    @(45375.99999999999,68,-73.50499495848095,208.4754760230378),
    @(45382.99999999999,67,-66.6248782391132,217.6680581371519),
    @(45389.99999999999,67,-74.9786672012971,214.4641185049499),
    @(45396.99999999999,67,-75.92929649062388,199.3220628226546),
    @(45403.99999999999,67,-84.48814186990654,215.6149832292396),
    @(45410.99999999999,67,-77.54041619597342,201.5596384264312)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

# --- 5) Match formatting to the rest of the workbook ---
# Header style (bold, centered, bordered) mirrors the other sheets' header row.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Date column style (custom date/time number format) mirrors column A elsewhere.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A32").PasteSpecial(-4122)

[void]($excel.CutCopyMode = $false)

# --- 6) Restore selection/active sheet like before the edit ---
[void]$wsForecast.Range("A1").Select()
$wsWeekly.Activate()
